# Weekly update: a new price-report row for "Apio" at Macroferia Regional de
# Talca needs to be inserted as the new row 147. All the existing rows that
# were at 147..178 shift down by one (to 148..179); their data is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 147, pushing rows 147:178 down to 148:179.
$ws.Rows("147:147").Insert()

# Populate the newly inserted row 147 with the new weekly record.
$ws.Range("A147").Value = 5
$ws.Range("B147").Value = 'Macroferia Regional de Talca'
$ws.Range("C147").Value = 'Maule'
$ws.Range("D147").Value = 44711
$ws.Range("E147").Value = 7
$ws.Range("F147").Value = 100112017
$ws.Range("G147").Value = 'Apio'
$ws.Range("H147").Value = 'Americana (o)'
$ws.Range("I147").Value = 'Primera'
$ws.Range("J147").Value = 700
$ws.Range("K147").Value = 6000
$ws.Range("L147").Value = 6000
$ws.Range("M147").Value = 6000
$ws.Range("N147").Value = '$/docena de matas'
$ws.Range("O147").Value = 'Provincia del Elquí'
$ws.Range("P147").Value = 1000
$ws.Range("Q147").Value = 6
$ws.Range("R147").Value = 'Hortaliza'
